$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Pais")

# Update the "last updated" timestamp string in A1
$ws.Range("A1").Value = "Datos actualizados a 21 de Abril de 2020 a las 19:52"

# Row 4 (e.g. Estados Unidos)
$ws.Range("B4").Value = 803575
$ws.Range("C4").Value = 10816
$ws.Range("E4").Value = 684595
$ws.Range("G4").Value = 1149
$ws.Range("H4").Value = 43663

# Row 16
$ws.Range("B16").Value = 38205
$ws.Range("C16").Value = 1376
$ws.Range("E16").Value = 23788
$ws.Range("G16").Value = 141
$ws.Range("H16").Value = 1831

# Row 22
$ws.Range("B22").Value = 16040
$ws.Range("C22").Value = 388
$ws.Range("E22").Value = 15233
$ws.Range("F22").Value = 315
$ws.Range("G22").Value = 43
$ws.Range("H22").Value = 730

# Row 75
$ws.Range("E75").Value = 1014
$ws.Range("G75").Value = 1
$ws.Range("H75").Value = 38

# Row 145
$ws.Range("D145").Value = 26
$ws.Range("E145").Value = 65
$ws.Range("F145").Value = 2
